$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste")

# New "indsatspakker" rows appended after the existing list (row 83 was last).
# Shared-string table order matters: populate new unique strings in the
# same order the author's workbook recorded them.
$ws.Range("A85").Value = "Dag - opstartsindsatser SEL § 83 a"

$ws.Range("A84").Value = "Aften - opstartsindsatser SEL § 83 a"
$ws.Range("A84").NumberFormat = "@"

$ws.Range("A86").Value = "Dag -- Terminalpakke Servicelov"
$ws.Range("A86").NumberFormat = "@"

$ws.Range("A87").Value = "Aften -- Terminalpakke Servicelov"
$ws.Range("A87").NumberFormat = "@"

$ws.Range("A87").Select()
